# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values (currentAveragePrice* / LevePrice* / LeveProfit*)
# to the Omega_Profits workbook, one hunk/row at a time, matching the upstream diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 4329103.5
$ws.Range("I9").Value = 6493556
$ws.Range("K9").Value = 6493556
$ws.Range("M9").Value = -6493387

$ws.Range("H12").Value = 367.61765
$ws.Range("I12").Value = 299.96875
$ws.Range("K12").Value = 299.96875
$ws.Range("M12").Value = -129.96875

$ws.Range("H32").Value = 3269.6924
$ws.Range("I32").Value = 951.7
$ws.Range("K32").Value = 951.7
$ws.Range("M32").Value = -625.7

$ws.Range("H40").Value = 2983.2856
$ws.Range("I40").Value = 2321
$ws.Range("J40").Value = 3866.3333
$ws.Range("K40").Value = 2321
$ws.Range("L40").Value = 3866.3333
$ws.Range("M40").Value = -2146
$ws.Range("N40").Value = -4216.3333

$ws.Range("H87").Value = 249995
$ws.Range("J87").Value = 249995
$ws.Range("L87").Value = 249995
$ws.Range("N87").Value = -252491

$ws.Range("H90").Value = 249995
$ws.Range("J90").Value = 249995
$ws.Range("L90").Value = 749985
$ws.Range("N90").Value = -762465

$ws.Range("H132").Value = 3726.2666
$ws.Range("I132").Value = 3473.76
$ws.Range("K132").Value = 10421.28
$ws.Range("M132").Value = -7891.280000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 281.875
$ws.Range("I4").Value = 307.75
$ws.Range("J4").Value = 256
$ws.Range("K4").Value = 307.75
$ws.Range("L4").Value = 256
$ws.Range("M4").Value = -191.75
$ws.Range("N4").Value = -488

$ws.Range("H5").Value = 157.3
$ws.Range("I5").Value = 141.44444
$ws.Range("K5").Value = 141.44444
$ws.Range("M5").Value = -29.44443999999999

$ws.Range("H9").Value = 20000
$ws.Range("I9").Value = 20000
$ws.Range("K9").Value = 20000
$ws.Range("M9").Value = -19830

$ws.Range("H20").Value = 20000
$ws.Range("I20").Value = 20000
$ws.Range("K20").Value = 20000
$ws.Range("M20").Value = -19730

$ws.Range("H46").Value = 19332.334
$ws.Range("J46").Value = 18998
$ws.Range("L46").Value = 18998
$ws.Range("N46").Value = -19636

$ws.Range("H61").Value = 8741.412
$ws.Range("I61").Value = 5468.4165
$ws.Range("J61").Value = 16596.6
$ws.Range("K61").Value = 5468.4165
$ws.Range("L61").Value = 16596.6
$ws.Range("M61").Value = -5256.4165
$ws.Range("N61").Value = -17020.6

$ws.Range("H74").Value = 2392.35
$ws.Range("I74").Value = 2392.35
$ws.Range("K74").Value = 2392.35
$ws.Range("M74").Value = -1518.35

$ws.Range("H77").Value = 2392.35
$ws.Range("I77").Value = 2392.35
$ws.Range("K77").Value = 11961.75
$ws.Range("M77").Value = -7593.75

$ws.Range("H132").Value = 3275.0967
$ws.Range("I132").Value = 3173.724
$ws.Range("K132").Value = 9521.172
$ws.Range("M132").Value = -6991.172

$ws.Range("H136").Value = 8741.412
$ws.Range("I136").Value = 5468.4165
$ws.Range("J136").Value = 16596.6
$ws.Range("K136").Value = 16405.2495
$ws.Range("L136").Value = 49789.8
$ws.Range("M136").Value = -13855.2495
$ws.Range("N136").Value = -54889.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 157.3
$ws.Range("I4").Value = 141.44444
$ws.Range("K4").Value = 141.44444
$ws.Range("M4").Value = -26.44443999999999

$ws.Range("H20").Value = 3514.1428
$ws.Range("I20").Value = 3520
$ws.Range("J20").Value = 3499.5
$ws.Range("K20").Value = 3520
$ws.Range("L20").Value = 3499.5
$ws.Range("M20").Value = -3273
$ws.Range("N20").Value = -3993.5

$ws.Range("H86").Value = 8335456
$ws.Range("I86").Value = 10418764
$ws.Range("J86").Value = 2224.5
$ws.Range("K86").Value = 10418764
$ws.Range("L86").Value = 2224.5
$ws.Range("M86").Value = -10417641
$ws.Range("N86").Value = -4470.5

$ws.Range("H88").Value = 24999.5
$ws.Range("J88").Value = 24999.5
$ws.Range("L88").Value = 24999.5
$ws.Range("N88").Value = -25811.5

$ws.Range("H89").Value = 8335456
$ws.Range("I89").Value = 10418764
$ws.Range("J89").Value = 2224.5
$ws.Range("K89").Value = 52093820
$ws.Range("L89").Value = 11122.5
$ws.Range("M89").Value = -52088204
$ws.Range("N89").Value = -22354.5

$ws.Range("H91").Value = 24999.5
$ws.Range("J91").Value = 24999.5
$ws.Range("L91").Value = 24999.5
$ws.Range("N91").Value = -27807.5

$ws.Range("H134").Value = 3389.6858
$ws.Range("I134").Value = 3339.4119
$ws.Range("K134").Value = 10018.2357
$ws.Range("M134").Value = -7483.235700000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 60.923077
$ws.Range("I7").Value = 51.8
$ws.Range("K7").Value = 51.8
$ws.Range("M7").Value = 61.2

$ws.Range("H22").Value = 1000.5833
$ws.Range("I22").Value = 852.8
$ws.Range("J22").Value = 1106.1428
$ws.Range("K22").Value = 852.8
$ws.Range("L22").Value = 1106.1428
$ws.Range("M22").Value = -502.8
$ws.Range("N22").Value = -1806.1428

$ws.Range("H88").Value = 34299.8
$ws.Range("I88").Value = 31311
$ws.Range("J88").Value = 35047
$ws.Range("K88").Value = 31311
$ws.Range("L88").Value = 35047
$ws.Range("N88").Value = -35859
$ws.Range("M88").Value = -30905

$ws.Range("H91").Value = 34299.8
$ws.Range("I91").Value = 31311
$ws.Range("J91").Value = 35047
$ws.Range("K91").Value = 31311
$ws.Range("L91").Value = 35047
$ws.Range("N91").Value = -37855
$ws.Range("M91").Value = -29907

$ws.Range("H93").Value = 20500
$ws.Range("I93").Value = 20500
$ws.Range("K93").Value = 20500
$ws.Range("M93").Value = -18628

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1049.5
$ws.Range("J68").Value = 1049.5
$ws.Range("L68").Value = 3148.5
$ws.Range("N68").Value = -4770.5

$ws.Range("H71").Value = 1049.5
$ws.Range("J71").Value = 1049.5
$ws.Range("L71").Value = 9445.5
$ws.Range("N71").Value = -17557.5

$ws.Range("H134").Value = 9499.546
$ws.Range("I134").Value = 1066
$ws.Range("J134").Value = 12662.125
$ws.Range("K134").Value = 3198
$ws.Range("L134").Value = 37986.375
$ws.Range("M134").Value = 1872
$ws.Range("N134").Value = -48126.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 175.1875
$ws.Range("I2").Value = 137.23077
$ws.Range("K2").Value = 137.23077
$ws.Range("M2").Value = -24.23077000000001

$ws.Range("H80").Value = 4900.0586
$ws.Range("I80").Value = 3095.8333
$ws.Range("K80").Value = 3095.8333
$ws.Range("M80").Value = -2097.8333

$ws.Range("H83").Value = 4900.0586
$ws.Range("I83").Value = 3095.8333
$ws.Range("K83").Value = 15479.1665
$ws.Range("M83").Value = -10487.1665

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

$ws.Range("H113").Value = 6145.5386
$ws.Range("J113").Value = 3999
$ws.Range("L113").Value = 3999
$ws.Range("N113").Value = -8339

$ws.Range("H132").Value = 3472.2764
$ws.Range("I132").Value = 3488.2878
$ws.Range("J132").Value = 3366.6
$ws.Range("K132").Value = 10464.8634
$ws.Range("L132").Value = 10099.8
$ws.Range("M132").Value = -7934.8634
$ws.Range("N132").Value = -15159.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1442.2858
$ws.Range("I16").Value = 1068.8889
$ws.Range("J16").Value = 3682.6667
$ws.Range("K16").Value = 1068.8889
$ws.Range("L16").Value = 3682.6667
$ws.Range("M16").Value = -898.8888999999999
$ws.Range("N16").Value = -4022.6667

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H25").Value = 12000
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H100").Value = 3136.25
$ws.Range("I100").Value = 1975.2
$ws.Range("J100").Value = 3664
$ws.Range("K100").Value = 1975.2
$ws.Range("L100").Value = 3664
$ws.Range("M100").Value = -1434.2
$ws.Range("N100").Value = -4746

$ws.Range("H132").Value = 4557
$ws.Range("I132").Value = 3379.8
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 10139.4
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -7609.400000000001
$ws.Range("N132").Value = -27560

$ws.Range("H138").Value = 89999
$ws.Range("J138").Value = 89999
$ws.Range("L138").Value = 89999
$ws.Range("N138").Value = -100279

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2283
$ws.Range("I81").Value = 2147.15
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 4294.3
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -3233.3
$ws.Range("N81").Value = -12122

$ws.Range("H84").Value = 2283
$ws.Range("I84").Value = 2147.15
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 21471.5
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -16167.5
$ws.Range("N84").Value = -60608

$ws.Range("H132").Value = 11308.875
$ws.Range("I132").Value = 10079.667
$ws.Range("J132").Value = 14996.5
$ws.Range("K132").Value = 30239.001
$ws.Range("L132").Value = 44989.5
$ws.Range("M132").Value = -27709.001
$ws.Range("N132").Value = -50049.5

$ws.Range("H136").Value = 3549.4167
$ws.Range("I136").Value = 3610.3215
$ws.Range("K136").Value = 10830.9645
$ws.Range("M136").Value = -8280.9645
